$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValue = 53.71147335634279

$ws.Range("N2").Value = $newValue
$ws.Range("N3").Value = $newValue
$ws.Range("N4").Value = $newValue
$ws.Range("N5").Value = $newValue
$ws.Range("N6").Value = $newValue
